# 2022-Q1 data was added to 603233-大参林.xlsx.
#
# Concretely (matching the authored diff):
#   1. The existing "总计" (grand-total) sheet is renamed to "2022-Q1" and its
#      contents are replaced by the new quarter's per-fund holding detail.
#   2. A brand-new "总计" sheet is created right after "2022-Q1", carrying the
#      same summary table as before plus one new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# A sheet that already contains the bold/centered/bordered "header & index
# column" style used throughout this workbook (style used by B1 and A2).
$styleSource = $wb.Worksheets.Item("2021-Q4")

function Copy-HeaderStyle($ws, $row, $col) {
    $styleSource.Range("B1").Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

function Set-TextCell($ws, $row, $col, $text) {
    # Force a General-formatted cell to hold a literal text value even when
    # the text looks numeric (fund codes with leading zeros, "25.68", ...).
    # A leading apostrophe is Excel's own "treat as text" quote-prefix; the
    # stored value itself does not include the apostrophe.
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$detail = $wb.Worksheets.Item("总计")
$detail.Name = "2022-Q1"

$detailHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $detailHeaders.Length; $c++) {
    $col = $c + 2
    $detail.Cells.Item(1, $col).Value = $detailHeaders[$c]
    Copy-HeaderStyle $detail 1 $col
}

$detailRows = @(
    @("002408", "中信建投医改灵活配置混合A", "25.68", "94.92", "4.66", "1.1967", 7),
    @("007553", "中信建投医改灵活配置混合C", "13.34", "94.92", "4.66", "0.6216", 7),
    @("010090", "中信建投医药健康混合A",     "5.37",  "94.87", "4.77", "0.2561", 8),
    @("010091", "中信建投医药健康混合C",     "2.25",  "94.87", "4.77", "0.1073", 8),
    @("501007", "汇添富中证互联网医疗主题指数（LOF）A", "0.58", "93.89", "4.67", "0.0271", 10),
    @("005043", "国寿安保健康科学混合A",     "0.99",  "85.72", "2.72", "0.0269", 8),
    @("005044", "国寿安保健康科学混合C",     "0.87",  "85.72", "2.72", "0.0237", 8),
    @("009502", "国寿安保创新医药股票A",     "0.54",  "81.60", "2.79", "0.0151", 8),
    @("501008", "汇添富中证互联网医疗主题指数（LOF）C", "0.19", "93.89", "4.67", "0.0089", 10),
    @("009503", "国寿安保创新医药股票C",     "0.20",  "81.60", "2.79", "0.0056", 8)
)

for ($i = 0; $i -lt $detailRows.Length; $i++) {
    $r = $i + 2
    $row = $detailRows[$i]

    $detail.Cells.Item($r, 1).Value = $i
    Copy-HeaderStyle $detail $r 1

    Set-TextCell $detail $r 2 $row[0]
    Set-TextCell $detail $r 3 $row[1]
    Set-TextCell $detail $r 4 $row[2]
    Set-TextCell $detail $r 5 $row[3]
    Set-TextCell $detail $r 6 $row[4]
    Set-TextCell $detail $r 7 $row[5]
    $detail.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: create the new "总计" sheet right after "2022-Q1"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $detail)
$total.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($c = 0; $c -lt $totalHeaders.Length; $c++) {
    $col = $c + 2
    $total.Cells.Item(1, $col).Value = $totalHeaders[$c]
    Copy-HeaderStyle $total 1 $col
}

$totalRows = @(
    @("2022-Q1", 10, 2.29),
    @("2021-Q4", 19, 5.05),
    @("2021-Q3", 12, 2.97),
    @("2021-Q2", 21, 5.96),
    @("2021-Q1", 51, 16.58),
    @("2020-Q4", 41, 20.39)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    $total.Cells.Item($r, 1).Value = $i
    Copy-HeaderStyle $total $r 1

    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}
